# repull data, push all data, mean calculation
# Update the dSF column (F) values for several rows to reflect re-pulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -7
    4  = 1
    6  = -2
    11 = -7
    12 = -3
    14 = -1
    17 = -2
    18 = 5
    19 = -2
    20 = 6
    21 = -3
    23 = 8
    24 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
